# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for rows 2-25 on the active
# worksheet, matching the recalculated s_vals from the regenerated save data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value
$kValues = @{
    2  = 0
    3  = 3
    4  = 0
    5  = 2
    6  = 3
    7  = 3
    8  = 3
    9  = 2
    10 = 7
    11 = 2
    12 = 6
    13 = 4
    14 = 3
    15 = 4
    16 = 2
    17 = 3
    18 = 2
    19 = 4
    20 = 8
    21 = 3
    22 = 4
    23 = 3
    24 = 1
    25 = 4
}

foreach ($row in $kValues.Keys | Sort-Object) {
    $ws.Range("G$row").Value = $kValues[$row]
}
